$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38 (shifts old rows 38..57 down to 39..58)
$ws.Rows.Item(38).Insert()

# Populate the new row 38 with the weekly price-report entry
$ws.Cells.Item(38, 1).Value = 11
$ws.Cells.Item(38, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(38, 3).Value = "Bíobío"
$ws.Cells.Item(38, 4).Value = 44518
$ws.Cells.Item(38, 5).Value = 8
$ws.Cells.Item(38, 6).Value = 100112021
$ws.Cells.Item(38, 7).Value = "Ají"
$ws.Cells.Item(38, 8).Value = "Inferno"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 50
$ws.Cells.Item(38, 11).Value = 14000
$ws.Cells.Item(38, 12).Value = 15000
$ws.Cells.Item(38, 13).Value = 14400
$ws.Cells.Item(38, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(38, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(38, 16).Value = 1200
$ws.Cells.Item(38, 17).Value = 12
$ws.Cells.Item(38, 18).Value = "Hortaliza"
